$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Replace the retailer data (rows 2-8) with the new list of retailers
# ---------------------------------------------------------------------------
$data = @(
  @("RET-08803", "Poroshi Telecom", "Sardah Baza"),
  @("RET-21070", "Tajmul Telecom", "Station Market"),
  @("RET-23822", "Mohona Telecom", "Sherkul Bazar"),
  @("RET-26501", "Tamanna Telecom", "Tebaria Baza"),
  @("RET-26506", "Tripty Electronics", "Jonail bazar"),
  @("RET-29332", "Rasel Electronics", "Bagha Bazar"),
  @("RET-32048", "Islam Enterprise", "Singra")
)

$row = 2
foreach ($item in $data) {
  $ws.Cells.Item($row, 1).Value = $item[0]
  $ws.Cells.Item($row, 2).Value = $item[1]
  $ws.Cells.Item($row, 3).Value = $item[2]
  $row++
}

# Rows 9 and 10 no longer hold any retailer data - clear their contents
# (formatting of those two rows is left untouched)
$ws.Range("A9:C10").ClearContents()

# ---------------------------------------------------------------------------
# 2. Re-align the header row and the data table so text is centered both
#    horizontally and vertically
# ---------------------------------------------------------------------------

# Header row (A1:C1) - bold header cells, center horizontally + vertically
$ws.Range("A1:C1").HorizontalAlignment = -4108
$ws.Range("A1:C1").VerticalAlignment = -4108

# RetailerID column cells (A2 and A5) only get centered (no wrap)
$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A2").VerticalAlignment = -4108
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)

# Remaining data cells (RetailerName / RetailerAddress columns, plus the
# RetailerID cells in rows 3,4,6,7,8) get centered + wrapped
$ws.Range("B2").HorizontalAlignment = -4108
$ws.Range("B2").VerticalAlignment = -4108
$ws.Range("B2").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("B3:C8").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("A3:A4").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("A6:A8").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Add the small red note cell in D4 (mirrors the existing style used for
#    the helper cells in column N)
# ---------------------------------------------------------------------------
$ws.Range("N6").Copy()
$ws.Range("D4").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4. Update the active cell selection
# ---------------------------------------------------------------------------
$ws.Range("F10").Select()
